# Auto-generated script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.434.23"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "3.767.96"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "3.765.94"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.64%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000253"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "4.397.09"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "3.772.04"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "69.469.45"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.119"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000135"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  -4.60%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "460.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.62%  "
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "2.954.15"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "27.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.38%  "

Write-Output "Applied 85 cell updates"
